# Add 4 new header rows at the top of the sheet for contact info,
# pushing the existing Table1 data down (A1:O2 -> A5:O6), and apply
# alternating yellow / theme-accent fills with bold labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 rows above the existing table (shifts table + everything down)
$ws.Range("A1:A4").EntireRow.Insert()

# Row 1: Company Name:
$ws.Range("A1").Value = "Company Name:"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 65535  # yellow (FFFF00 BGR)
$ws.Range("B1").Interior.ThemeColor = 8
$ws.Range("B1").Interior.TintAndShade = 0.39997558519241921

# Row 2: Your Name:
$ws.Range("A2").Value = "Your Name:"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Interior.ThemeColor = 8
$ws.Range("A2").Interior.TintAndShade = 0.39997558519241921
$ws.Range("B2").Interior.Color = 65535

# Row 3: Phone Number:
$ws.Range("A3").Value = "Phone Number:"
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Interior.Color = 65535
$ws.Range("B3").Interior.ThemeColor = 8
$ws.Range("B3").Interior.TintAndShade = 0.39997558519241921

# Row 4: Email ID:
$ws.Range("A4").Value = "Email ID:"
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Interior.ThemeColor = 8
$ws.Range("A4").Interior.TintAndShade = 0.39997558519241921
$ws.Range("B4").Interior.Color = 65535

# Update selection to match target
$ws.Range("B13").Select()
